# Commit: "add alert when user withdraw required course"
#
# The function sheet (功能表) tracks which feature requirements are
# implemented (O) vs not yet implemented (X). Row 11 corresponds to
# "退選必修課需提出警告" (an alert must be raised when a student withdraws
# from a required/compulsory course). This feature has now been
# implemented, so its status cell flips from "X" to "O".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("功能表")

# Flip the status flag for "退選必修課需提出警告" from X to O.
$ws.Range("B11").Value = "O"

# Reflect that B11 is the cell the author was just working on.
$ws.Range("B11").Select() | Out-Null
